$wb = $excel.ActiveWorkbook

# Helper: write a value as TEXT (shared string) without leaving a residual
# NumberFormat style on the cell -- mirrors how this workbook's data-entry
# app stores numeric-looking values (phone numbers, quantities, balances)
# as plain text rather than numbers.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# ---- Sheet "Stocks": reduce remaining quantities (col D) for items sold ----
$stocks = $wb.Worksheets.Item("Stocks")
$stocks.Cells.Item(2,4).Value = 201
$stocks.Cells.Item(4,4).Value = 54
$stocks.Cells.Item(6,4).Value = 98
$stocks.Cells.Item(8,4).Value = 0
$stocks.Cells.Item(9,4).Value = 148
$stocks.Cells.Item(11,4).Value = 41

# ---- Sheet "Bills": clear balance on two existing rows (Bill Clearance) ----
$bills = $wb.Worksheets.Item("Bills")
Set-TextValue $bills.Cells.Item(3,5) "0"
Set-TextValue $bills.Cells.Item(6,5) "0"

# ---- Sheet "Bills": append new bill log rows 7-22 ----
# Row 7
$bills.Cells.Item(7,1).Value = "02-Dec-2020 14:09"
$bills.Cells.Item(7,2).Value = "Muthu Rathinam"
Set-TextValue $bills.Cells.Item(7,3) "66559988"
Set-TextValue $bills.Cells.Item(7,4) "1000"
Set-TextValue $bills.Cells.Item(7,5) "0"
$bills.Cells.Item(7,6).Value = "XX021284"
$bills.Cells.Item(7,7).Value = "Stock Name Three(1)"

# Row 8
$bills.Cells.Item(8,1).Value = "02-Dec-2020 14:09"
$bills.Cells.Item(8,2).Value = "Muthu Rathinam"
Set-TextValue $bills.Cells.Item(8,3) "66559988"
Set-TextValue $bills.Cells.Item(8,4) "175"
Set-TextValue $bills.Cells.Item(8,5) "0"
$bills.Cells.Item(8,6).Value = "XX021285"
$bills.Cells.Item(8,7).Value = "cutting(1)"

# Row 9
$bills.Cells.Item(9,1).Value = "02-Dec-2020 14:11"
$bills.Cells.Item(9,2).Value = "Muthu Rathinam"
Set-TextValue $bills.Cells.Item(9,3) "66559988"
Set-TextValue $bills.Cells.Item(9,4) "35"
Set-TextValue $bills.Cells.Item(9,5) "0"
$bills.Cells.Item(9,6).Value = "XX021285"
$bills.Cells.Item(9,7).Value = "Sak 2(1)"

# Row 10
$bills.Cells.Item(10,1).Value = "02-Dec-2020 14:22"
$bills.Cells.Item(10,2).Value = "Muthu Rathinam"
Set-TextValue $bills.Cells.Item(10,3) "66559988"
Set-TextValue $bills.Cells.Item(10,4) "35"
Set-TextValue $bills.Cells.Item(10,5) "0"
$bills.Cells.Item(10,6).Value = "XX021286"
$bills.Cells.Item(10,7).Value = "sakthi masala(1)"

# Row 11
$bills.Cells.Item(11,1).Value = "02-Dec-2020 14:25"
$bills.Cells.Item(11,2).Value = "Muthu Rathinam"
Set-TextValue $bills.Cells.Item(11,3) "66559988"
Set-TextValue $bills.Cells.Item(11,4) "35"
Set-TextValue $bills.Cells.Item(11,5) "0"
$bills.Cells.Item(11,6).Value = "XX021286"
$bills.Cells.Item(11,7).Value = "sakthi masala(1)"

# Row 12
$bills.Cells.Item(12,1).Value = "02-Dec-2020 14:26"
$bills.Cells.Item(12,2).Value = "Muthu Rathinam"
Set-TextValue $bills.Cells.Item(12,3) "66559988"
Set-TextValue $bills.Cells.Item(12,4) "685"
Set-TextValue $bills.Cells.Item(12,5) "0"
$bills.Cells.Item(12,6).Value = "XX021286"
$bills.Cells.Item(12,7).Value = "sakthi masala(1),Bill Clearance 02Dec2020(1)"

# Row 13
$bills.Cells.Item(13,1).Value = "02-Dec-2020 14:32"
$bills.Cells.Item(13,2).Value = "Muthu Rathinam"
Set-TextValue $bills.Cells.Item(13,3) "66559988"
Set-TextValue $bills.Cells.Item(13,4) "1600"
Set-TextValue $bills.Cells.Item(13,5) "0"
$bills.Cells.Item(13,6).Value = "XX021287"
$bills.Cells.Item(13,7).Value = "CP(2)"

# Row 14
$bills.Cells.Item(14,1).Value = "02-Dec-2020 14:32"
$bills.Cells.Item(14,2).Value = "Muthu Rathinam"
Set-TextValue $bills.Cells.Item(14,3) "66559988"
Set-TextValue $bills.Cells.Item(14,4) "2250"
Set-TextValue $bills.Cells.Item(14,5) "0.0"
$bills.Cells.Item(14,6).Value = "XX021287"
$bills.Cells.Item(14,7).Value = "CP(2),Bill Clearance 02Dec2020(1)"

# Row 15
$bills.Cells.Item(15,1).Value = "02-Dec-2020 14:34"
$bills.Cells.Item(15,2).Value = "Aseth"
Set-TextValue $bills.Cells.Item(15,3) "88979"
Set-TextValue $bills.Cells.Item(15,4) "2400"
Set-TextValue $bills.Cells.Item(15,5) "0"
$bills.Cells.Item(15,6).Value = "XX021288"
$bills.Cells.Item(15,7).Value = "CP(3)"

# Row 16
$bills.Cells.Item(16,1).Value = "02-Dec-2020 14:34"
$bills.Cells.Item(16,2).Value = "Aseth"
Set-TextValue $bills.Cells.Item(16,3) "88979"
Set-TextValue $bills.Cells.Item(16,4) "70"
Set-TextValue $bills.Cells.Item(16,5) "0"
$bills.Cells.Item(16,6).Value = "XX021289"
$bills.Cells.Item(16,7).Value = "Sak 2(2)"

# Row 17
$bills.Cells.Item(17,1).Value = "02-Dec-2020 14:35"
$bills.Cells.Item(17,2).Value = "Aseth"
Set-TextValue $bills.Cells.Item(17,3) "88979"
Set-TextValue $bills.Cells.Item(17,4) "175"
Set-TextValue $bills.Cells.Item(17,5) "0"
$bills.Cells.Item(17,6).Value = "XX021290"
$bills.Cells.Item(17,7).Value = "Sak 2(5)"

# Row 18
$bills.Cells.Item(18,1).Value = "02-Dec-2020 14:35"
$bills.Cells.Item(18,2).Value = "Aseth"
Set-TextValue $bills.Cells.Item(18,3) "88979"
Set-TextValue $bills.Cells.Item(18,4) "800"
Set-TextValue $bills.Cells.Item(18,5) "0"
$bills.Cells.Item(18,6).Value = "XX021291"
$bills.Cells.Item(18,7).Value = "CP(1)"

# Row 19
$bills.Cells.Item(19,1).Value = "02-Dec-2020 14:39"
$bills.Cells.Item(19,2).Value = "Aseth"
Set-TextValue $bills.Cells.Item(19,3) "88979"
Set-TextValue $bills.Cells.Item(19,4) "140"
Set-TextValue $bills.Cells.Item(19,5) "0"
$bills.Cells.Item(19,6).Value = "XX021291"
$bills.Cells.Item(19,7).Value = "Single Motta (4)"

# Row 20
$bills.Cells.Item(20,1).Value = "02-Dec-2020 14:39"
$bills.Cells.Item(20,2).Value = "Aseth"
Set-TextValue $bills.Cells.Item(20,3) "88979"
Set-TextValue $bills.Cells.Item(20,4) "35"
Set-TextValue $bills.Cells.Item(20,5) "0"
$bills.Cells.Item(20,6).Value = "XX021292"
$bills.Cells.Item(20,7).Value = "Sak 2(1)"

# Row 21
$bills.Cells.Item(21,1).Value = "02-Dec-2020 14:40"
$bills.Cells.Item(21,2).Value = "Aseth"
Set-TextValue $bills.Cells.Item(21,3) "88979"
Set-TextValue $bills.Cells.Item(21,4) "155"
Set-TextValue $bills.Cells.Item(21,5) "0"
$bills.Cells.Item(21,6).Value = "XX021293"
$bills.Cells.Item(21,7).Value = "Stock Name One(1)"

# Row 22
$bills.Cells.Item(22,1).Value = "02-Dec-2020 14:40"
$bills.Cells.Item(22,2).Value = "Aseth"
Set-TextValue $bills.Cells.Item(22,3) "88979"
Set-TextValue $bills.Cells.Item(22,4) "195"
Set-TextValue $bills.Cells.Item(22,5) "0.0"
$bills.Cells.Item(22,6).Value = "XX021293"
$bills.Cells.Item(22,7).Value = "Stock Name One(1),Bill Clearance 02Dec2020(1)"

Write-Host "Applied Bills log update (Yes / No Done)"
